$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The table's last row (58) was mistakenly dated "14/7/2026". A new
# weekly report for 7/16/2025 is being added. While doing so, the
# typo'd date on the existing last row is corrected to "15/7/2025"
# (the row itself is unchanged otherwise), and the brand new row (59)
# becomes "16/7/2025" with its own progress numbers.

# Copy formatting (borders, fonts, number formats, etc.) from row 58
# down into the new row 59 before putting values into it.
$ws.Range("D58:J58").Copy() | Out-Null
$ws.Range("D59:J59").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item(59).RowHeight = $ws.Rows.Item(58).RowHeight

# Fill in the new row first so the "16/7/2025" shared string is created
# before the existing row's text is rewritten to "15/7/2025".
$ws.Range("D59").Value = "16/7/2025"
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 924
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 1012
$ws.Range("J59").Value = "N/A"

# Now correct the old row's date text.
$ws.Range("D58").Value = "15/7/2025"

# Grow the Excel table (ListObject) so the new row is included.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("D4:J59"))

# Keep the view roughly where the diff shows it ended up.
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("F63").Select()
